$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: shared string "1" -> "2". A plain Value assignment of a numeric-looking
# string gets auto-converted to a number by Excel, so force text entry via
# NumberFormat, then restore the original (General/style-0) formatting by
# copying formats back from an untouched cell.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "2"
$ws.Range("C1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# B1: shared string "stri" -> "laptop acer" (already non-numeric text, no
# auto-conversion risk).
$ws.Range("B1").Value = "laptop acer"

# C1: numeric 10.0 -> 100.0
$ws.Range("C1").Value = 100.0

# D1: numeric 12.0 -> 10.0
$ws.Range("D1").Value = 10.0
